$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with two new columns P and Q, copying formatting from O1
$ws.Range("O1").Copy($ws.Range("P1:Q1"))
$ws.Range("P1").Value2 = 14
$ws.Range("Q1").Value2 = 15

# For each data row (2 to 25): swap values between columns I/K and M/O, then add P and Q = 2
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2   # column I = 9
    $kVal = $ws.Cells.Item($r, 11).Value2  # column K = 11
    $mVal = $ws.Cells.Item($r, 13).Value2  # column M = 13
    $oVal = $ws.Cells.Item($r, 15).Value2  # column O = 15

    $ws.Cells.Item($r, 9).Value2 = $kVal
    $ws.Cells.Item($r, 11).Value2 = $iVal
    $ws.Cells.Item($r, 13).Value2 = $oVal
    $ws.Cells.Item($r, 15).Value2 = $mVal

    $ws.Cells.Item($r, 16).Value2 = 2  # column P = 16
    $ws.Cells.Item($r, 17).Value2 = 2  # column Q = 17
}
